# Generate Report for Handoff
# Adds a new entry for aa856427-b3d2-4f7e-a140-e7d1e57eef81.md to the
# Overview / zh-cn / de-de localization-status tables.

$wb = $excel.ActiveWorkbook

$fileName      = "aa856427-b3d2-4f7e-a140-e7d1e57eef81.md"
$pathAndName   = "e2e\aa856427-b3d2-4f7e-a140-e7d1e57eef81.md"
$ext           = ".md"
$status        = "Ready for handoff"
$sourcePath    = "e2e"
$priority      = "ht"
$contentDup    = "False"
$zhXlfName     = "aa856427-b3d2-4f7e-a140-e7d1e57eef81.b23bc585a7d799e32d310b11649693554e5bcd0e.zh-cn.xlf"
$deXlfName     = "aa856427-b3d2-4f7e-a140-e7d1e57eef81.b23bc585a7d799e32d310b11649693554e5bcd0e.de-de.xlf"
$handoffDate   = "2016-10-18 12:11:29"
$zhHandoffDt   = "2016-10-18 12:11:16"
$deHandoffDt   = "2016-10-18 12:11:29"
$handbackDt    = "0001-01-01 00:00:00"
$refTokens     = "True"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa856427d4e3f0c1a2b3c4d5e6f708192a3b4c5d/e2e/aa856427-b3d2-4f7e-a140-e7d1e57eef81.md"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Range("A8").Value = $fileName
$wsOverview.Range("B8").Value = $pathAndName
$wsOverview.Range("C8").Value = $ext
$wsOverview.Range("D8").Value = ""
$wsOverview.Range("E8").Value = $status
$wsOverview.Range("F8").Value = $status
$wsOverview.Range("G8").Value = $handoffDate
$wsOverview.Range("G8").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B8"), $ghBase, "", "", $pathAndName)

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()

$wsZh.Range("A8").Value = $fileName
$wsZh.Range("B8").Value = $ext
$wsZh.Range("C8").Value = $status
$wsZh.Range("D8").Value = $sourcePath
$wsZh.Range("E8").Value = $priority
$wsZh.Range("F8").Value = $contentDup
$wsZh.Range("G8").Value = $zhXlfName
$wsZh.Range("H8").Value = $zhHandoffDt
$wsZh.Range("H8").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I8").Value = ""
$wsZh.Range("J8").Value = ""
$wsZh.Range("K8").Value = $handbackDt
$wsZh.Range("K8").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L8").Value = ""
$wsZh.Range("M8").Value = $refTokens
$wsZh.Range("N8").Value = ""
$wsZh.Range("O8").Value = $contentDup
$wsZh.Range("P8").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A8"), $ghBase, "", "", $fileName)

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()

$wsDe.Range("A8").Value = $fileName
$wsDe.Range("B8").Value = $ext
$wsDe.Range("C8").Value = $status
$wsDe.Range("D8").Value = $sourcePath
$wsDe.Range("E8").Value = $priority
$wsDe.Range("F8").Value = $contentDup
$wsDe.Range("G8").Value = $deXlfName
$wsDe.Range("H8").Value = $deHandoffDt
$wsDe.Range("H8").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I8").Value = ""
$wsDe.Range("J8").Value = ""
$wsDe.Range("K8").Value = $handbackDt
$wsDe.Range("K8").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L8").Value = ""
$wsDe.Range("M8").Value = $refTokens
$wsDe.Range("N8").Value = ""
$wsDe.Range("O8").Value = $contentDup
$wsDe.Range("P8").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A8"), $ghBase, "", "", $fileName)

Write-Host "Done adding handback row for $fileName"
